# "Add ability to load external files into TPA" commit — the part of the
# change that touches this workbook is cosmetic: rename the sheet and move
# the saved selection down one row below the last used row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Sheet1 -> All
$ws.Name = "All"

# Selection moves from C19 to A20 (first empty row under the table)
$ws.Range("A20").Select()
